# "I have changed last page"
#
# The author appended one brand-new slide to the very end of the deck.
# The new slide uses the same layout as slide 1 ("TITLE" / Title Slide,
# i.e. slideLayout1.xml) and is left essentially blank: an empty title
# placeholder ("Title 1") and an empty subtitle placeholder
# ("Subtitle 2"). Nothing on any of the other slides changes.

$p = $ppt.ActivePresentation

$slideCount = $p.Slides.Count

# The Title Slide custom layout is the one used by slide 1 (matches the
# placeholder set seen in the new slide: ctrTitle + subTitle).
$master = $p.SlideMaster
$titleLayout = $master.CustomLayouts.Item(1)

# Append a brand-new slide after the current last slide.
$newSlide = $p.Slides.AddSlide($slideCount + 1, $titleLayout)

# Give the two placeholders the same shape names PowerPoint assigns to a
# freshly-inserted, still-empty title slide.
$newSlide.Shapes.Item(1).Name = "Title 1"
$newSlide.Shapes.Item(2).Name = "Subtitle 2"
